$d = $word.ActiveDocument

# The document contains two paragraphs with the text
# "The Result Activity will look like this:" — only the FIRST one
# (the caption right after the Hunt Activity picture, i.e. the one that
# precedes the Result Activity screenshot) is renamed to
# "The Start Activity will look like this:" per the commit.
$targetParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*The Result Activity will look like this:*") {
        $targetParaIndex = $i
        break
    }
}

if ($targetParaIndex -eq -1) {
    throw "Could not locate target paragraph"
}

$p = $d.Paragraphs.Item($targetParaIndex)
$paraStart = $p.Range.Start
$paraText = $p.Range.Text

$wordIdx = $paraText.IndexOf("Result")
$resultStart = $paraStart + $wordIdx
$resultEnd = $resultStart + "Result".Length

# Replace "Result" with "Start" (typed-replacement, like selecting the
# word and retyping it).
$resultRange = $d.Range($resultStart, $resultEnd)
$resultRange.Text = "Start"

# Re-apply the explicit character formatting to the newly typed word so
# it becomes its own run (matching the split seen in the authored edit:
# "The " / "Start" / " Activity will look like this:").
$newWordEnd = $resultStart + "Start".Length
$newWordRange = $d.Range($resultStart, $newWordEnd)
$newWordRange.Font.Name = "Times New Roman"
$newWordRange.Font.Size = 12

Write-Output "Paragraph $targetParaIndex now reads: $($d.Paragraphs.Item($targetParaIndex).Range.Text)"
